$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

$ws.Range("C2").Value = 245
$ws.Range("D2").Value = 99.19028340080972
$ws.Range("C3").Value = 245
$ws.Range("D3").Value = 99.19028340080972
$ws.Range("C4").Value = 162
$ws.Range("D4").Value = 65.58704453441295
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2.42914979757085
$ws.Range("C7").Value = 42
$ws.Range("D7").Value = 17.00404858299595
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 4.8582995951417
$ws.Range("C9").Value = 228
$ws.Range("D9").Value = 92.30769230769231
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0.8097165991902834
$ws.Range("C12").Value = 58
$ws.Range("D12").Value = 23.48178137651822
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 12.14574898785425
$ws.Range("C15").Value = 158
$ws.Range("D15").Value = 63.96761133603239
$ws.Range("C16").Value = 97
$ws.Range("D16").Value = 39.27125506072874
$ws.Range("C18").Value = 109
$ws.Range("D18").Value = 44.12955465587044
$ws.Range("C19").Value = 218
$ws.Range("D19").Value = 88.25910931174089
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 5.668016194331984
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 1.214574898785425
$ws.Range("C22").Value = 214
$ws.Range("D22").Value = 86.63967611336032
$ws.Range("C23").Value = 214
$ws.Range("D23").Value = 86.63967611336032
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 6.477732793522267
$ws.Range("C25").Value = 41
$ws.Range("D25").Value = 16.59919028340081
$ws.Range("C26").Value = 60
$ws.Range("D26").Value = 24.2914979757085
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 5.668016194331984
$ws.Range("C29").Value = 34
$ws.Range("D29").Value = 13.76518218623482
$ws.Range("C30").Value = 45
$ws.Range("D30").Value = 18.21862348178137
$ws.Range("C31").Value = 227
$ws.Range("D31").Value = 91.90283400809717
$ws.Range("C32").Value = 160
$ws.Range("D32").Value = 64.77732793522267
$ws.Range("C33").Value = 18
$ws.Range("D33").Value = 7.28744939271255
$ws.Range("C34").Value = 212
$ws.Range("D34").Value = 85.82995951417004
$ws.Range("C35").Value = 97
$ws.Range("D35").Value = 39.27125506072874
